$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (Strike#) column values regenerated from s_vals calc
$kValues = @{
    2  = 5
    3  = 0
    4  = 4
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    9  = 2
    10 = 3
    11 = 1
    12 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
